# Week 17 data logging + running totals update for the Giants 2021 Team Data
# workbook, across the YDS / OFF / DEF / ST / TURNS / PEN sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# YDS sheet: append this week's per-play yardage log to each running string.
# ---------------------------------------------------------------------------
$ydsWs = $wb.Worksheets.Item("YDS")

$ydsWs.Range("B2").Value2 = $ydsWs.Range("B2").Value2 + " 8 4 6 2 6 5 4 -1 8 5 1 2 10 8 2 3 3 5 2 -2 4 -3 -1 9 6 4 1 8 3 3 1 1 13 5 7 5 1 3 9 1"
$ydsWs.Range("C2").Value2 = $ydsWs.Range("C2").Value2 + " 2 0 -1 4 10 5 3 4 7 2 7 7 9 2 2 2 10 4 4 1 4 -1 1 3 3 0"
$ydsWs.Range("B3").Value2 = $ydsWs.Range("B3").Value2 + " 4 12 12 -4"
$ydsWs.Range("C3").Value2 = $ydsWs.Range("C3").Value2 + " 3 11 8 4 6 23 18 17 4 6 16 8 8 11 10 3 13"

# ---------------------------------------------------------------------------
# OFF sheet: season running totals after Week 17.
# ---------------------------------------------------------------------------
$offWs = $wb.Worksheets.Item("OFF")

$offWs.Range("C2").Value2 = 179
$offWs.Range("E2").Value2 = 9
$offWs.Range("F2").Value2 = 78
$offWs.Range("G2").Value2 = 46
$offWs.Range("I2").Value2 = 8
$offWs.Range("J2").Value2 = 27
$offWs.Range("N2").Value2 = 20
$offWs.Range("O2").Value2 = 26

$offWs.Range("C3").Value2 = 199
$offWs.Range("G3").Value2 = 32
$offWs.Range("I3").Value2 = 76
$offWs.Range("J3").Value2 = 58
$offWs.Range("L3").Value2 = 316
$offWs.Range("M3").Value2 = 183
$offWs.Range("Q3").Value2 = 561

# ---------------------------------------------------------------------------
# DEF sheet: season running totals after Week 17.
# ---------------------------------------------------------------------------
$defWs = $wb.Worksheets.Item("DEF")

$defWs.Range("C2").Value2 = 212
$defWs.Range("D2").Value2 = 13
$defWs.Range("F2").Value2 = 67
$defWs.Range("G2").Value2 = 72
$defWs.Range("I2").Value2 = 7
$defWs.Range("J2").Value2 = 32
$defWs.Range("N2").Value2 = 13
$defWs.Range("O2").Value2 = 22
$defWs.Range("P2").Value2 = 15

$defWs.Range("B3").Value2 = 13
$defWs.Range("C3").Value2 = 192
$defWs.Range("E3").Value2 = 28
$defWs.Range("F3").Value2 = 113
$defWs.Range("G3").Value2 = 38
$defWs.Range("H3").Value2 = 19
$defWs.Range("I3").Value2 = 64
$defWs.Range("J3").Value2 = 65
$defWs.Range("L3").Value2 = 334
$defWs.Range("M3").Value2 = 224
$defWs.Range("Q3").Value2 = 626

# ---------------------------------------------------------------------------
# ST sheet: kick/return logs and season totals.
# ---------------------------------------------------------------------------
$stWs = $wb.Worksheets.Item("ST")

$stWs.Range("B2").Value2 = 63
$stWs.Range("D2").Value2 = 68
$stWs.Range("J2").Value2 = 177
$stWs.Range("K2").Value2 = 160
$stWs.Range("B3").Value2 = 27

$stWs.Range("B4").Value2 = $stWs.Range("B4").Value2 + " 53 51"
$stWs.Range("B5").Value2 = $stWs.Range("B5").Value2 + " 16 6"
$stWs.Range("B6").Value2 = $stWs.Range("B6").Value2 + " 26 24 3"
$stWs.Range("D3").Value2 = $stWs.Range("D3").Value2 + " 46 39 33 38"
$stWs.Range("D4").Value2 = $stWs.Range("D4").Value2 + " 1 0 0 0"
$stWs.Range("D5").Value2 = $stWs.Range("D5").Value2 + " 15 0 17 0 0"

# ---------------------------------------------------------------------------
# TURNS sheet: Road turnover totals.
# ---------------------------------------------------------------------------
$turnsWs = $wb.Worksheets.Item("TURNS")

$turnsWs.Range("B3").Value2 = 12
$turnsWs.Range("C3").Value2 = 6
$turnsWs.Range("D3").Value2 = 10
$turnsWs.Range("E3").Value2 = 7

# ---------------------------------------------------------------------------
# PEN sheet: penalty totals.
# ---------------------------------------------------------------------------
$penWs = $wb.Worksheets.Item("PEN")

$penWs.Range("B2").Value2 = 19
$penWs.Range("D4").Value2 = 5
